$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 890304.5600000001
$ws.Range("I9").Value = 230.58333
$ws.Range("J9").Value = 1602363.8
$ws.Range("K9").Value = 230.58333
$ws.Range("L9").Value = 1602363.8
$ws.Range("M9").Value = -61.58332999999999
$ws.Range("N9").Value = -1602701.8
$ws.Range("H88").Value = 3126.125
$ws.Range("I88").Value = 2433.3333
$ws.Range("J88").Value = 3541.8
$ws.Range("K88").Value = 2433.3333
$ws.Range("L88").Value = 3541.8
$ws.Range("M88").Value = -2027.3333
$ws.Range("N88").Value = -4353.8
$ws.Range("H91").Value = 3126.125
$ws.Range("I91").Value = 2433.3333
$ws.Range("J91").Value = 3541.8
$ws.Range("K91").Value = 2433.3333
$ws.Range("L91").Value = 3541.8
$ws.Range("M91").Value = -1029.3333
$ws.Range("N91").Value = -6349.8
$ws.Range("H137").Value = 7330.4
$ws.Range("I137").Value = 2581.182
$ws.Range("J137").Value = 13135
$ws.Range("K137").Value = 7743.545999999999
$ws.Range("L137").Value = 39405
$ws.Range("M137").Value = -5193.545999999999
$ws.Range("N137").Value = -44505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 500
$ws.Range("K4").Value = 500
$ws.Range("M4").Value = -384
$ws.Range("H5").Value = 963.3333
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("N5").Value = -1224
$ws.Range("H61").Value = 4284.3657
$ws.Range("I61").Value = 3958.2856
$ws.Range("K61").Value = 3958.2856
$ws.Range("M61").Value = -3746.2856
$ws.Range("H64").Value = 5030000
$ws.Range("I64").Value = 10000000
$ws.Range("K64").Value = 10000000
$ws.Range("M64").Value = -9999752
$ws.Range("H67").Value = 5030000
$ws.Range("I67").Value = 10000000
$ws.Range("K67").Value = 10000000
$ws.Range("M67").Value = -9999142
$ws.Range("H97").Value = 2315603.2
$ws.Range("J97").Value = 561.25
$ws.Range("L97").Value = 561.25
$ws.Range("N97").Value = -1553.25
$ws.Range("H109").Value = 41627.668
$ws.Range("J109").Value = 41627.668
$ws.Range("L109").Value = 41627.668
$ws.Range("N109").Value = -44401.668
$ws.Range("H110").Value = 11365433
$ws.Range("I110").Value = 20834128
$ws.Range("K110").Value = 20834128
$ws.Range("M110").Value = -20832083
$ws.Range("H136").Value = 4284.3657
$ws.Range("I136").Value = 3958.2856
$ws.Range("K136").Value = 11874.8568
$ws.Range("M136").Value = -9324.856800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 963.3333
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 1000
$ws.Range("N4").Value = -1230
$ws.Range("H25").Value = 6213.35
$ws.Range("I25").Value = 912.1667
$ws.Range("J25").Value = 8485.286
$ws.Range("K25").Value = 912.1667
$ws.Range("L25").Value = 8485.286
$ws.Range("M25").Value = -677.1667
$ws.Range("N25").Value = -8955.286
$ws.Range("H92").Value = 90000
$ws.Range("J92").Value = 90000
$ws.Range("L92").Value = 90000
$ws.Range("N92").Value = -94992
$ws.Range("H108").Value = 97975.2
$ws.Range("J108").Value = 97975.2
$ws.Range("L108").Value = 97975.2
$ws.Range("N108").Value = -105655.2
$ws.Range("H132").Value = 98075.75
$ws.Range("J132").Value = 106372.29
$ws.Range("L132").Value = 106372.29
$ws.Range("N132").Value = -116492.29

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 168.36363
$ws.Range("I7").Value = 217.125
$ws.Range("J7").Value = 38.333332
$ws.Range("K7").Value = 217.125
$ws.Range("L7").Value = 38.333332
$ws.Range("M7").Value = -104.125
$ws.Range("N7").Value = -264.333332
$ws.Range("H31").Value = 3783.4688
$ws.Range("I31").Value = 2092.0715
$ws.Range("K31").Value = 2092.0715
$ws.Range("M31").Value = -1797.0715
$ws.Range("H34").Value = 3783.4688
$ws.Range("I34").Value = 2092.0715
$ws.Range("K34").Value = 2092.0715
$ws.Range("M34").Value = -1890.0715
$ws.Range("H56").Value = 19449.5
$ws.Range("J56").Value = 19900
$ws.Range("L56").Value = 19900
$ws.Range("N56").Value = -21590
$ws.Range("H96").Value = 25999.666
$ws.Range("J96").Value = 25999.666
$ws.Range("L96").Value = 25999.666
$ws.Range("N96").Value = -31491.666
$ws.Range("H134").Value = 3125.1667
$ws.Range("J134").Value = 5628.909
$ws.Range("L134").Value = 16886.727
$ws.Range("N134").Value = -21956.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3236.8635
$ws.Range("I2").Value = 108.5
$ws.Range("J2").Value = 3932.0557
$ws.Range("K2").Value = 651
$ws.Range("L2").Value = 23592.3342
$ws.Range("M2").Value = -538
$ws.Range("N2").Value = -23818.3342
$ws.Range("H11").Value = 15118.223
$ws.Range("J11").Value = 6665.6665
$ws.Range("L11").Value = 19996.9995
$ws.Range("N11").Value = -20276.9995
$ws.Range("H38").Value = 42.75
$ws.Range("I38").Value = 27.833334
$ws.Range("K38").Value = 83.50000199999999
$ws.Range("M38").Value = 263.499998
$ws.Range("H55").Value = 1524.5
$ws.Range("J55").Value = 1766.1666
$ws.Range("L55").Value = 5298.4998
$ws.Range("N55").Value = -5652.4998
$ws.Range("H80").Value = 3387
$ws.Range("I80").Value = 1994.3334
$ws.Range("J80").Value = 4083.3333
$ws.Range("K80").Value = 5983.0002
$ws.Range("L80").Value = 12249.9999
$ws.Range("M80").Value = -5047.0002
$ws.Range("N80").Value = -14121.9999
$ws.Range("H83").Value = 3387
$ws.Range("I83").Value = 1994.3334
$ws.Range("J83").Value = 4083.3333
$ws.Range("K83").Value = 17949.0006
$ws.Range("L83").Value = 36749.9997
$ws.Range("M83").Value = -13269.0006
$ws.Range("N83").Value = -46109.9997
$ws.Range("H86").Value = 606
$ws.Range("J86").Value = 699.5
$ws.Range("L86").Value = 2098.5
$ws.Range("N86").Value = -4470.5
$ws.Range("H89").Value = 606
$ws.Range("J89").Value = 699.5
$ws.Range("L89").Value = 6295.5
$ws.Range("N89").Value = -18151.5
$ws.Range("H131").Value = 9012846
$ws.Range("I131").Value = 41667576
$ws.Range("J131").Value = 4644.8623
$ws.Range("K131").Value = 125002728
$ws.Range("L131").Value = 13934.5869
$ws.Range("M131").Value = -124997688
$ws.Range("N131").Value = -24014.5869

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 570.7143
$ws.Range("I2").Value = 750
$ws.Range("J2").Value = 122.5
$ws.Range("K2").Value = 750
$ws.Range("L2").Value = 122.5
$ws.Range("M2").Value = -637
$ws.Range("N2").Value = -348.5
$ws.Range("H119").Value = 70000
$ws.Range("J119").Value = 70000
$ws.Range("L119").Value = 70000
$ws.Range("N119").Value = -79676

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4312.5864
$ws.Range("I132").Value = 3641.2964
$ws.Range("K132").Value = 10923.8892
$ws.Range("M132").Value = -8393.889200000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 19500
$ws.Range("I10").Value = 19500
$ws.Range("K10").Value = 19500
$ws.Range("M10").Value = -19331
$ws.Range("H96").Value = 700
$ws.Range("I96").Value = 700
$ws.Range("K96").Value = 700
$ws.Range("M96").Value = 673
